$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.025.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.90%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.674.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.95%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'216.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.26%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.529"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.39%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +2.75%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.83%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.15%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +3.98%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.910.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.97%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.685.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.78%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.09%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +2.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'27.035.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.90%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'233.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.75%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0737"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.53%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.26%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.14%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +3.42%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.12%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'145.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.25%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.02%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.56%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'15.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.54%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.07%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.78%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.89%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.455.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.50%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +5.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +5.15%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'0.899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.51%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.83%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.66%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +3.64%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.12%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +3.69%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.980"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.71%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +5.18%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.816.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.84%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.784"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.93%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'90.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.82%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.28%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +4.28%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D51").Value = "'7.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.79%  "
$ws.Range("E51").Style = "Normal"
